# Trading update: 2026-02-17 19:47:32
# Append the newest MarketMaking trade (trade #10, still OPEN) as row 11
# on both the "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 11

    $ws.Cells.Item($row, 1).Value = 10

    # Column B ("2026-02-17") must land as literal text, not an auto-parsed
    # date serial: stage it as a formula that evaluates to the text, then
    # paste-special just the resulting value into place (keeps default style).
    $ws.Range("Z1").Formula = '="2026-02-17"'
    $ws.Range("Z1").Copy()
    $ws.Range("B" + $row).PasteSpecial(-4163)
    $ws.Range("Z1").ClearContents()

    $ws.Cells.Item($row, 3).Value = "19:47:27"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.83
    $ws.Range("G" + $row).Formula = '=""'
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.0027272727273
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Range("P" + $row).Formula = '=""'
    $ws.Cells.Item($row, 17).Value = 0
}
